$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(1, -0.46865769738307961, 0.46686837043186813),
    @(2, -0.3462172945012103, 0.34080332796068724),
    @(3, -0.13567242477264685, 0.13504293487991959),
    @(4, -0.12304293497257035, 0.12249876618871269),
    @(5, -0.11649876651384172, 0.11542960212888609),
    @(6, -0.048897964744414768, 0.048864771317711941),
    @(7, -0.028864771714145476, 0.02882399773812061),
    @(8, -0.0088239981360098696, 0.0088176437017892795),
    @(9, -0.0028176440372762457, 0.0028146776129940321),
    @(10, 0.0031853220513866631, -0.0031845476895355773),
    @(11, -0.051482554028503102, 0.051402318269545333),
    @(12, -0.045402318606510672, 0.04515409651364255),
    @(13, -0.039154096855932963, 0.039086490381000871),
    @(14, -0.027086490752003201, 0.027053779699781622),
    @(15, -0.021053780045035886, 0.02102803311566781),
    @(16, -0.015028033462132662, 0.015004530570141839),
    @(17, -0.0090045309181876476, 0.0089999996380356606),
    @(18, -0.11065107625352155, 0.11048623324072437),
    @(19, -0.027097430354861451, 0.027013861145249951),
    @(20, -0.018013861467926162, 0.018004299918702671),
    @(21, -0.0090043002418562779, 0.0089999996764866808),
    @(22, -0.19613090351959528, 0.19446293081965926),
    @(23, -0.13227716054181471, 0.13063827000284434),
    @(24, -0.088638270499173721, 0.087777373258987268),
    @(25, -0.060531637717151199, 0.060435502262201624),
    @(26, -0.054435502595655549, 0.054317338597655862),
    @(27, -0.048317338932292841, 0.047933171217785286),
    @(28, -0.041933171556684634, 0.041683160572477007),
    @(29, -0.02968316094106882, 0.029580442865359657),
    @(30, -0.0095804432710999698, 0.009560227957024825),
    @(31, -0.027019868363030497, 0.027000837734604488),
    @(32, -0.0060008381462060001, 0.0059999996557786872),
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
}
